$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove columns G:L entirely (data shrinks from A1:L4 to A1:F4)
$ws.Range("G1:L1").EntireColumn.Delete()

# Header row (row 1) - plain text values, no quoting needed
$ws.Range("A1").Value = "ad"
$ws.Range("B1").Value = "aw"
$ws.Range("C1").Value = "da"
$ws.Range("D1").Value = "dw"
$ws.Range("E1").Value = "wa"
$ws.Range("F1").Value = "wd"

# Data rows - values look numeric, so write with a leading apostrophe to
# force them to remain text cells like the original inline strings, then
# reset the style so no quotePrefix formatting sticks on the cell.
$ws.Range("A2").Value = "'0.795"
$ws.Range("B2").Value = "'0.825"
$ws.Range("C2").Value = "'0.567"
$ws.Range("D2").Value = "'0.926"
$ws.Range("E2").Value = "'0.607"
$ws.Range("F2").Value = "'0.982"

$ws.Range("A3").Value = "'0.843"
$ws.Range("B3").Value = "'0.815"
$ws.Range("C3").Value = "'0.587"
$ws.Range("D3").Value = "'0.897"
$ws.Range("E3").Value = "'0.628"
$ws.Range("F3").Value = "'0.974"

$ws.Range("A4").Value = "'0.833"
$ws.Range("B4").Value = "'0.824"
$ws.Range("C4").Value = "'0.586"
$ws.Range("D4").Value = "'0.913"
$ws.Range("E4").Value = "'0.612"
$ws.Range("F4").Value = "'0.970"

# Strip the quote-prefix style the apostrophe entry adds so cells have no
# explicit style index, matching the source workbook's plain formatting.
$ws.Range("A2:F4").Style = "Normal"
